$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Dditems"
$ws.Range("B32").Value = "present"
$ws.Range("C32").Value = "absent"
$ws.Range("D32").Value = "late"
$ws.Range("E32").Value = "excused"

$ws.Range("A32").Select()
